$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Formula = "'246.38"

# Row 3
$ws.Cells.Item(3, 4).Formula = "'22.28"

# Row 4
$ws.Cells.Item(4, 4).Formula = "'5.355"

# Row 5
$ws.Cells.Item(5, 4).Formula = "'0.05861"

# Row 7
$ws.Cells.Item(7, 4).Formula = "'6.376"

# Row 8
$ws.Cells.Item(8, 4).Formula = "'0.8135"

# Row 9
$ws.Cells.Item(9, 4).Formula = "'0.9996"

# Row 10
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Formula = "'0.1423"
$ws.Cells.Item(10, 5).Value = "9WazirXWRX"

# Row 11
$ws.Cells.Item(11, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(11, 4).Formula = "'0.03873"
$ws.Cells.Item(11, 5).Value = "10LiechtensteinCryptoassetsExchangeLCX"

# Row 12
$ws.Cells.Item(12, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12, 4).Formula = "'0.07328"
$ws.Cells.Item(12, 5).Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).Formula = "'0.03005"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

# Row 14
$ws.Cells.Item(14, 2).Value = "MCDex"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(14, 4).Formula = "'4.181"
$ws.Cells.Item(14, 5).Value = "13MCDexMCB"

# Row 15
$ws.Cells.Item(15, 2).Value = "BitMartToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(15, 4).Formula = "'0.09394"
$ws.Cells.Item(15, 5).Value = "14BitMartTokenBMX"

# Row 16
$ws.Cells.Item(16, 2).Value = "BitForexToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(16, 4).Formula = "'0.001593"
$ws.Cells.Item(16, 5).Value = "15BitForexTokenBF"

# Row 17
$ws.Cells.Item(17, 2).Value = "CoinExToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(17, 4).Formula = "'0.04823"
$ws.Cells.Item(17, 5).Value = "16CoinExTokenCET"

# Row 18
$ws.Cells.Item(18, 2).Value = "One"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(18, 4).Formula = "'0.0005890"
$ws.Cells.Item(18, 5).Value = "17OneONE"

# Row 19
$ws.Cells.Item(19, 4).Formula = "'0.006039"

# Row 20
$ws.Cells.Item(20, 4).Formula = "'0.004081"

# Row 21
$ws.Cells.Item(21, 4).Formula = "'0.0009834"

# Row 22
$ws.Cells.Item(22, 4).Formula = "'0.0001300"

# Row 23
$ws.Cells.Item(23, 4).Formula = "'3.689"

# Row 25
$ws.Cells.Item(25, 4).Formula = "'0.3245"

# Row 40
$ws.Cells.Item(40, 4).Formula = "'0.03850"

# Row 41
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 4).Formula = "'0.1074"
$ws.Cells.Item(41, 5).Value = "40BKEXTokenBKK"

# Row 42
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42, 4).Formula = "'0.002410"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"

# Row 43
$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43, 4).Formula = "'0.003048"
$ws.Cells.Item(43, 5).Value = "42KickTokenKICKWorstin24h"

# Row 44
$ws.Cells.Item(44, 4).Formula = "'0.005205"

# Row 45
$ws.Cells.Item(45, 4).Formula = "'0.00005645"

# Row 47
$ws.Cells.Item(47, 4).Formula = "'0.7220"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOINBestin24h"

# Row 48
$ws.Cells.Item(48, 4).Formula = "'0.08546"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"

# Row 49
$ws.Cells.Item(49, 4).Formula = "'0.00002100"
